# Reword heading of agency overview section
# ---------------------------------------------------------------
# 1. "Goal status agency" -> "Agency overview" (document heading)
# 2. Collapse the split "{{" / "agency_name" / "}}" runs (and the
#    proofErr spell-check markers around them) into one run reading
#    "Goal Status across {{agency_name}}"
# 3. Collapse the split "{{" / "agency_abbreviation" / "}}" runs into
#    one run
# 4. Collapse the split "{{" / "goal_change_summary_sentence" / "}}"
#    runs into one run
# 5. Collapse the split "{{" / "challenge_summary_text" / "}}" runs
#    into one run
# 6. Drop the proofErr spell-check markers that wrap
#    "goal_status_breakdown_bullets" while keeping its run (and its
#    "{{" / "r " / "}}" neighbours) separate, as in the source edit.
# ---------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Heading text -------------------------------------------------
$d.Content.Find.Execute("Goal status agency", $true, $false, $false, $false, $false, $true, 1, $false, "Agency overview", 2) | Out-Null

# --- 2. "Goal Status across {{agency_name}}" --------------------------
$d.Content.Find.Execute("Goal Status across {{agency_name}}", $true, $false, $false, $false, $false, $true, 1, $false, "Goal Status across {{agency_name}}", 2) | Out-Null

# --- 3. "{{agency_abbreviation}}" --------------------------------------
$d.Content.Find.Execute("{{agency_abbreviation}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{agency_abbreviation}}", 2) | Out-Null

# --- 4. "{{goal_change_summary_sentence}}" -----------------------------
$d.Content.Find.Execute("{{goal_change_summary_sentence}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{goal_change_summary_sentence}}", 2) | Out-Null

# --- 5. "{{challenge_summary_text}}" ------------------------------------
$d.Content.Find.Execute("{{challenge_summary_text}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{challenge_summary_text}}", 2) | Out-Null

# --- 6. Strip proofErr around "goal_status_breakdown_bullets" ----------
# Find the paragraph that holds the "{{r goal_status_breakdown_bullets}}"
# placeholder so we can scope all further operations to just it.
$count = $d.Paragraphs.Count
$bulletsParaIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "goal_status_breakdown_bullets") {
        $bulletsParaIndex = $i
    }
}

if ($bulletsParaIndex -gt 0) {
    $bp = $d.Paragraphs.Item($bulletsParaIndex)
    $pStart = $bp.Range.Start
    $pEnd = $bp.Range.End

    # Re-write the whole placeholder text in place: this regenerates the
    # run (merging the four fragments into one, formatting unchanged)
    # and drops the stale spell-check proofErr markers in the process.
    $mergeRng = $d.Range($pStart, $pEnd)
    $mergeRng.Find.Execute("{{r goal_status_breakdown_bullets}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{r goal_status_breakdown_bullets}}", 2) | Out-Null

    # Re-split that merged run back into the original four text runs
    # ("{{", "r ", "goal_status_breakdown_bullets", "}}") by toggling a
    # character property on/off across each exact sub-range -- this
    # creates run boundaries without reintroducing proofErr markers.
    $segments = @(
        @(0, 2),   # {{
        @(2, 4),   # "r "
        @(4, 33),  # goal_status_breakdown_bullets
        @(33, 35)  # }}
    )
    foreach ($seg in $segments) {
        $s = $pStart + $seg[0]
        $e = $pStart + $seg[1]
        $onRng = $d.Range($s, $e)
        $onRng.Bold = 1
        $offRng = $d.Range($s, $e)
        $offRng.Bold = 0
    }
}
